$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0000000002517112612985156
$ws.Range("B3").Value = 0.000000007679762822577009
$ws.Range("B4").Value = 0.00000009273040197581414
$ws.Range("B5").Value = 0.0000006725351380400409
$ws.Range("B6").Value = 0.000003521222073310612
$ws.Range("B7").Value = 0.00001467498748185719
$ws.Range("B8").Value = 0.00005161355800702858
$ws.Range("B9").Value = 0.0001591403836677027
$ws.Range("B10").Value = 0.0004416326008299711
$ws.Range("B11").Value = 0.001124308218147197
